$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "相談件数" - the only tab selected in the source file

# --- Add the new daily data row (2020-04-30) ---
# Copy the last existing data row (95) and insert a copy of it above row 96.
# This both (a) reuses the exact same number formats/styles as the rest of
# the table for the new row, and (b) naturally pushes the trailing footnote
# row (the "※4/8より..." note, originally on row 96) down to row 97.
$ws.Rows("95:95").Copy() | Out-Null
$ws.Rows("96:96").Insert(-4121) | Out-Null   # -4121 = xlShiftDown

# Overwrite the copied values with the real data for the new row
$ws.Range("A96").Value = 43951   # 2020-04-30
$ws.Range("B96").Value = 792
$ws.Range("C96").Value = 32302
$ws.Range("D96").Value = 185
$ws.Range("E96").Value = 6849

# --- Update the sheet's print area so it covers the newly added row ---
$printAreaName = $null
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $candidate = $wb.Names.Item($i)
    if ($candidate.Name -like "*Print_Area*") {
        $printAreaName = $candidate
    }
}
if ($printAreaName -ne $null) {
    $printAreaName.RefersTo = "=相談件数!`$A`$1:`$E`$99"
}

# --- Refresh the frozen-pane view so the newly visible row/column are selected ---
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("B2").Select() | Out-Null
$win.FreezePanes = $true
$ws.Range("D97").Select() | Out-Null
